$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.959.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.506'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.65'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.865.72'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.640.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.545'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0764'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.01'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.016.94'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.02'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.92'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.30'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.71'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.130'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.86'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.60'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.25'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.56%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.53'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.136.03'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.21'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.799'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.775.40'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.76'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.64'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.59%  '
